$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.974.80"
Set-TextValue "E2" "  +0.46%  "

Set-TextValue "D3" "1.642.03"
Set-TextValue "E3" "  +0.06%  "

Set-TextValue "E4" "  +0.35%  "

Set-TextValue "D5" "212.80"
Set-TextValue "E5" "  +0.37%  "

Set-TextValue "E6" "  +0.16%  "

Set-TextValue "E7" "  +0.40%  "

Set-TextValue "D8" "23.55"
Set-TextValue "E8" "  +0.76%  "

Set-TextValue "E9" "  -1.82%  "

Set-TextValue "E10" "  +0.39%  "

Set-TextValue "E11" "  +2.29%  "

Set-TextValue "D12" "1.875.36"
Set-TextValue "E12" "  +0.25%  "

Set-TextValue "D13" "1.651.15"
Set-TextValue "E13" "  +0.73%  "

Set-TextValue "E14" "  +0.71%  "

Set-TextValue "E15" "  +1.28%  "

Set-TextValue "D16" "65.58"
Set-TextValue "E16" "  +0.11%  "

Set-TextValue "D17" "27.977.68"
Set-TextValue "E17" "  +0.68%  "

Set-TextValue "D18" "233.28"
Set-TextValue "E18" "  +0.35%  "

Set-TextValue "E19" "  +0.36%  "

Set-TextValue "D20" "7.60"
Set-TextValue "E20" "  -0.63%  "

Set-TextValue "E21" "  +0.18%  "

Set-TextValue "D22" "10.57"
Set-TextValue "E22" "  -1.37%  "

Set-TextValue "E23" "  -0.41%  "

Set-TextValue "E24" "  -3.37%  "

Set-TextValue "D25" "152.83"
Set-TextValue "E25" "  +1.81%  "

Set-TextValue "E26" "  +0.05%  "

Set-TextValue "E27" "  -0.06%  "

Set-TextValue "E28" "  -0.10%  "

Set-TextValue "E29" "  +0.37%  "

Set-TextValue "E30" "  +0.83%  "

Set-TextValue "E31" "  +0.46%  "

Set-TextValue "E32" "  +3.38%  "

Set-TextValue "E33" "  +0.08%  "

Set-TextValue "D34" "1.409.94"
Set-TextValue "E34" "  -4.11%  "

Set-TextValue "E35" "  +1.56%  "

Set-TextValue "E36" "  +1.54%  "

Set-TextValue "B37" "VeChain"
Set-TextValue "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.0169"
Set-TextValue "E37" "  +1.01%  "

Set-TextValue "B38" "ImmutableX"
Set-TextValue "C38" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "0.565"
Set-TextValue "E38" "  +1.19%  "

Set-TextValue "B39" "TrustWalletToken"
Set-TextValue "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D39" "0.929"
Set-TextValue "E39" "  -0.55%  "

Set-TextValue "B40" "ARBITRUM"
Set-TextValue "C40" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "0.880"
Set-TextValue "E40" "  -0.55%  "

Set-TextValue "E41" "  +0.94%  "

Set-TextValue "E42" "  +0.19%  "

Set-TextValue "E43" "  +6.24%  "

Set-TextValue "D44" "67.28"
Set-TextValue "E44" "  -2.74%  "

Set-TextValue "D45" "5.53"
Set-TextValue "E45" "  +2.99%  "

Set-TextValue "E46" "  -2.92%  "

Set-TextValue "D47" "1.784.39"
Set-TextValue "E47" "  -0.04%  "

Set-TextValue "D48" "88.13"
Set-TextValue "E48" "  +0.33%  "

Set-TextValue "E49" "  -0.02%  "

Set-TextValue "E50" "  +0.43%  "

Set-TextValue "D51" "7.59"
Set-TextValue "E51" "  -1.30%  "
